$d = $word.ActiveDocument

# ============================================================
# PART A (bottom of doc first, so paragraph numbers above are
# unaffected): remove "map" table, insert "order" and "review"
# tables in its place.
# ============================================================

# The blank paragraph + "map" table occupy paragraphs 25-29 (1-indexed).
$startPara = $d.Paragraphs(25)
$endPara = $d.Paragraphs(29)
$delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$delRange.Delete()

# After the delete, paragraph 24 is product's closing "');'" and is now last in doc.
$anchorIdx = 24

$newLines = @(
 "create table order (",
 "p_no int not null,",
 "o_no int not null primary key auto_increment,",
 "o_time datetime not null DEFAULR CURRENT_TIMESTAMP,",
 "amount int not null,",
 "status int not null",
 ");",
 "",
 "create table review (",
 "u_no int not null,",
 "s_no int not null,",
 "content varchar(300) not null,",
 "image varchar(1000),",
 "r_time datetime not null DEFAULR CURRENT_TIMESTAMP",
 ");"
)

foreach ($line in $newLines) {
    $p = $d.Paragraphs($anchorIdx)
    $p.Range.InsertParagraphAfter()
    $anchorIdx = $anchorIdx + 1
    $d.Paragraphs($anchorIdx).Range.Text = $line
}

# Bold the table names "order" and "review"
$pOrderHeader = $d.Paragraphs(25)
$rOrderHeader = $pOrderHeader.Range
$rOrderHeader.Find.Execute("order ")
$rOrderHeader.Bold = 1

$pReviewHeader = $d.Paragraphs(33)
$rReviewHeader = $pReviewHeader.Range
$rReviewHeader.Find.Execute("review ")
$rReviewHeader.Bold = 1

# Move the _GoBack bookmark to right after "primary key" in the o_no line (paragraph 27)
$pONo = $d.Paragraphs(27)
$rONo = $pONo.Range
$rONo.Find.Execute("primary key")
$bookmarkPos = $rONo.End
$collapsed = $d.Range($bookmarkPos, $bookmarkPos)
$collapsed.Bookmarks.Add("_GoBack")

# ============================================================
# PART B: product table - rename its "u_no" column to "s_no"
# (paragraph 18, unaffected by Part A's edits further down).
# ============================================================
$p18 = $d.Paragraphs(18)
$r18 = $p18.Range
$r18.Find.Execute("u_no int not null,", $true, $false, $false, $false, $false, $true, 1, $false, "s_no int not null,", 2)

# ============================================================
# PART C: store table restructuring - split the "u_no" row into
# a new auto_increment "s_no" primary key plus a plain "u_no" row,
# add a new "name" column, and keep lat/lng/category but make
# lat/lng "not null". Net effect: 2 new rows in the store table.
# ============================================================

# paragraph 10: 'u_no not null primary key,' -> 's_no int not null auto_increment primary key,'
$p10 = $d.Paragraphs(10)
$r10 = $p10.Range
$r10.Find.Execute("u_no not null primary key,", $true, $false, $false, $false, $false, $true, 1, $false, "s_no int not null auto_increment primary key,", 2)

# insert new paragraph after 10 with "u_no int not null,"
$d.Paragraphs(10).Range.InsertParagraphAfter()
$d.Paragraphs(11).Range.Text = "u_no int not null,"

# paragraph 12: 'lat double,' -> 'name varchar(45) not null,'
$p12 = $d.Paragraphs(12)
$r12 = $p12.Range
$r12.Find.Execute("lat double,", $true, $false, $false, $false, $false, $true, 1, $false, "name varchar(45) not null,", 2)

# paragraph 13: 'lng double,' -> 'lat double not null,'
$p13 = $d.Paragraphs(13)
$r13 = $p13.Range
$r13.Find.Execute("lng double,", $true, $false, $false, $false, $false, $true, 1, $false, "lat double not null,", 2)

# paragraph 14: 'category varchar(100) not null' -> 'lng double not null,'
$p14 = $d.Paragraphs(14)
$r14 = $p14.Range
$r14.Find.Execute("category varchar(100) not null", $true, $false, $false, $false, $false, $true, 1, $false, "lng double not null,", 2)

# paragraph 15: ');' -> 'category varchar(100) not null'
$p15 = $d.Paragraphs(15)
$r15 = $p15.Range
$r15.Find.Execute(");", $true, $false, $false, $false, $false, $true, 1, $false, "category varchar(100) not null", 2)

# insert new paragraph after 15 with ");"
$d.Paragraphs(15).Range.InsertParagraphAfter()
$d.Paragraphs(16).Range.Text = ");"

Write-Host "Final paragraph count: " $d.Paragraphs.Count
for ($i=1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Host $i ": [" $d.Paragraphs($i).Range.Text "]"
}
